$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.098.60"
$ws.Range("E2").Value = "  +0.10%  "

$ws.Range("D3").Value = "2.046.65"
$ws.Range("E3").Value = "  -0.69%  "

$ws.Range("E4").Value = "  +0.36%  "

$ws.Range("D5").Value = "249.11"
$ws.Range("E5").Value = "  -0.10%  "

$ws.Range("D6").Value = "0.667"
$ws.Range("E6").Value = "  -0.80%  "

$ws.Range("D7").Value = "59.41"
$ws.Range("E7").Value = "  +6.41%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("D9").Value = "0.386"
$ws.Range("E9").Value = "  +1.91%  "

$ws.Range("D10").Value = "0.0790"
$ws.Range("E10").Value = "  -1.33%  "

$ws.Range("E11").Value = "  +1.47%  "

$ws.Range("D12").Value = "15.99"
$ws.Range("E12").Value = "  +5.63%  "

$ws.Range("D13").Value = "2.345.98"
$ws.Range("E13").Value = "  -0.52%  "

$ws.Range("D14").Value = "0.834"
$ws.Range("E14").Value = "  +2.58%  "

$ws.Range("D15").Value = "5.79"
$ws.Range("E15").Value = "  +9.27%  "

$ws.Range("D16").Value = "2.044.05"
$ws.Range("E16").Value = "  -0.72%  "

$ws.Range("D17").Value = "18.50"
$ws.Range("E17").Value = "  +29.78%  "

$ws.Range("D18").Value = "37.061.24"
$ws.Range("E18").Value = "  +0.18%  "

$ws.Range("D19").Value = "75.62"
$ws.Range("E19").Value = "  +1.84%  "

$ws.Range("D20").Value = "0.0₃0902"
$ws.Range("E20").Value = "  -1.87%  "

$ws.Range("D21").Value = "5.37"
$ws.Range("E21").Value = "  -0.03%  "

$ws.Range("D22").Value = "237.66"
$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("E23").Value = "  -0.01%  "

$ws.Range("D24").Value = "2.42"
$ws.Range("E24").Value = "  +0.13%  "

$ws.Range("D25").Value = "2.19"
$ws.Range("E25").Value = "  +9.92%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "9.47"
$ws.Range("E26").Value = "  +4.19%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "169.58"
$ws.Range("E27").Value = "  -1.04%  "

$ws.Range("D28").Value = "20.16"
$ws.Range("E28").Value = "  +0.07%  "

$ws.Range("E29").Value = "  +0.94%  "

$ws.Range("D30").Value = "1.15"
$ws.Range("E30").Value = "  +8.08%  "

$ws.Range("D31").Value = "4.80"
$ws.Range("E31").Value = "  +4.59%  "

$ws.Range("D32").Value = "0.0631"
$ws.Range("E32").Value = "  +0.63%  "

$ws.Range("D33").Value = "4.56"
$ws.Range("E33").Value = "  +3.80%  "

$ws.Range("D34").Value = "0.0896"
$ws.Range("E34").Value = "  +3.41%  "

$ws.Range("E35").Value = "  +0.16%  "

$ws.Range("D36").Value = "2.23"
$ws.Range("E36").Value = "  -1.96%  "

$ws.Range("E37").Value = "  -1.54%  "

$ws.Range("D38").Value = "0.109"
$ws.Range("E38").Value = "  +3.77%  "

$ws.Range("E39").Value = "  -0.32%  "

$ws.Range("B40").Value = "THORChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D40").Value = "5.36"
$ws.Range("E40").Value = "  +22.65%  "

$ws.Range("B41").Value = "HuobiToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").Value = "3.22"
$ws.Range("E41").Value = "  +15.92%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "0.0224"
$ws.Range("E42").Value = "  -0.24%  "

$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "17.63"
$ws.Range("E43").Value = "  -3.06%  "

$ws.Range("E44").Value = "  -0.53%  "

$ws.Range("D45").Value = "97.01"
$ws.Range("E45").Value = "  +0.24%  "

$ws.Range("D46").Value = "2.50"
$ws.Range("E46").Value = "  +2.63%  "

$ws.Range("D47").Value = "1.295.47"
$ws.Range("E47").Value = "  -0.15%  "

$ws.Range("D49").Value = "6.88"
$ws.Range("E49").Value = "  +0.54%  "

$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.231.86"
$ws.Range("E50").Value = "  -0.77%  "

$ws.Range("B51").Value = "FTXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D51").Value = "3.53"
$ws.Range("E51").Value = "  -26.04%  "
